# Ajout du coût des plaquettes et DHT
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Plaquettes (PCB boards) - real cost now known
$ws.Range("G6").Value = 3

# Capteur de Temp et Humidité (DHT11) - replace placeholder note with real price calc
$ws.Range("F12").Value = $null
$ws.Range("E12").Value = 225
$ws.Range("F12").Value = 100
$ws.Range("G12").Formula = "=E12/F12"

# Update selection to match the author's final cursor position
$ws.Range("E14").Select()
